# calibration_wheel.xlsx — add the extra ordinary ray row-block for
# polarimetric acquisitions, fix the mistyped title, and restyle the
# header to match the data font.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the sheet title typo ("calibration wheell" -> "Calibration Wheel").
#    A1:D1 is merged, so writing the top-left cell is enough.
$ws.Range("A1").Value = "Calibration Wheel"

# 2) Re-style the header row (A1:D1) so it uses the same Times New Roman 12pt
#    font as the data rows, while keeping its own center/center alignment
#    and hair border. Paste the data-row format (A2) onto the header first —
#    this correctly re-uses the existing "Times New Roman" font resource
#    instead of registering a duplicate — then restore the header's vertical
#    centering that the paste overwrote.
$ws.Range("A2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4108
$excel.CutCopyMode = $false

# Header row is a touch taller in the new layout.
$ws.Rows.Item(1).RowHeight = 15

# 3) Widen column B slightly to fit the new "extra ordinary ray" labels.
$ws.Columns.Item(2).ColumnWidth = 16

# 4) Add 13 new blank rows (6-18) in columns A:B for the extra ordinary-ray
#    polarimetric acquisition entries. They share the data rows' font and
#    center/bottom alignment but carry no border.
$ws.Range("A2").Copy()
$ws.Range("A6:B18").PasteSpecial(-4122)
$ws.Range("A6:B18").Borders.LineStyle = -4142
$ws.Range("A6:B18").ClearContents()
$ws.Rows.Item("6:18").RowHeight = 15
$excel.CutCopyMode = $false

# 5) Restore the "home" selection so the file opens on A1 with A1:D5
#    highlighted, matching the reviewed layout.
$ws.Range("A1:D5").Select()
